$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing OrderDate values (E4, E7, E9)
$ws.Range("E4").Value = 42864
$ws.Range("E7").Value = 43870
$ws.Range("E9").Value = 44634

# Add new row 11 data
$ws.Range("A11").Value = "Tran"
$ws.Range("B11").Value = "Annette"
$ws.Range("C11").Value = 41000000022
$ws.Range("D11").Value = "annette_cigarette@gmail.com"

# Add hyperlink for D11, then reapply the same named "Hyperlink" cell style
# used by the existing rows so the cell's style index matches (rather than
# the fresh xf that Hyperlinks.Add creates on its own).
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:annette_cigarette@gmail.com")
$ws.Range("D11").Style = "Hyperlink"

# E11 date value, formatted like the other OrderDate cells
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = 44696

# Update column D width to fit new, longer content (closest reachable value
# to the recorded best-fit width of 26.1796875 characters)
$ws.Columns("D").ColumnWidth = 25.3

# Update selection to match recorded UI state
$ws.Range("I12").Select()
